{"js": "// Replace each \"two-digit \u00d7 two-digit = product\" answer in the table with\n// the new value from the commit. Every left-hand string below occurs\n// exactly once in the document, so an exact (non-wildcard) body.search()\n// safely targets the single matching run without touching anything else\n// (e.g. the unrelated date text at the top of the document).\nconst replacements = [\n  [\"74\u00d767=4958\", \"84\u00d775=6300\"],\n  [\"21\u00d726=546\", \"48\u00d724=1152\"],\n  [\"93\u00d733=3069\", \"45\u00d733=1485\"],\n  [\"66\u00d788=5808\", \"88\u00d717=1496\"],\n  [\"18\u00d717=306\", \"89\u00d798=8722\"],\n  [\"80\u00d799=7920\", \"89\u00d764=5696\"],\n  [\"62\u00d785=5270\", \"36\u00d730=1080\"],\n  [\"85\u00d781=6885\", \"94\u00d722=2068\"],\n  [\"23\u00d715=345\", \"77\u00d760=4620\"],\n  [\"49\u00d711=539\", \"89\u00d781=7209\"],\n  [\"96\u00d779=7584\", \"25\u00d774=1850\"],\n  [\"52\u00d785=4420\", \"21\u00d788=1848\"],\n  [\"95\u00d767=6365\", \"48\u00d782=3936\"],\n  [\"16\u00d725=400\", \"88\u00d757=5016\"],\n  [\"43\u00d756=2408\", \"46\u00d714=644\"],\n  [\"88\u00d731=2728\", \"62\u00d773=4526\"],\n  [\"42\u00d797=4074\", \"92\u00d767=6164\"],\n  [\"94\u00d749=4606\", \"19\u00d743=817\"],\n  [\"38\u00d721=798\", \"90\u00d753=4770\"],\n  [\"13\u00d780=1040\", \"64\u00d768=4352\"],\n  [\"80\u00d751=4080\", \"29\u00d793=2697\"],\n  [\"40\u00d716=640\", \"74\u00d745=3330\"],\n  [\"45\u00d728=1260\", \"92\u00d745=4140\"],\n  [\"58\u00d796=5568\", \"16\u00d712=192\"],\n  [\"27\u00d744=1188\", \"56\u00d796=5376\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"two-digit \u00d7 two-digit = product\" answer in the table with\n# the new value from the commit. Every find string below occurs exactly\n# once in the document, so Find.Execute with wdReplaceAll (applied to a\n# fresh whole-document range each time) safely updates only the single\n# matching run without touching anything else (e.g. the unrelated date\n# text at the top of the document).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"74\u00d767=4958\", \"84\u00d775=6300\"),\n    @(\"21\u00d726=546\", \"48\u00d724=1152\"),\n    @(\"93\u00d733=3069\", \"45\u00d733=1485\"),\n    @(\"66\u00d788=5808\", \"88\u00d717=1496\"),\n    @(\"18\u00d717=306\", \"89\u00d798=8722\"),\n    @(\"80\u00d799=7920\", \"89\u00d764=5696\"),\n    @(\"62\u00d785=5270\", \"36\u00d730=1080\"),\n    @(\"85\u00d781=6885\", \"94\u00d722=2068\"),\n    @(\"23\u00d715=345\", \"77\u00d760=4620\"),\n    @(\"49\u00d711=539\", \"89\u00d781=7209\"),\n    @(\"96\u00d779=7584\", \"25\u00d774=1850\"),\n    @(\"52\u00d785=4420\", \"21\u00d788=1848\"),\n    @(\"95\u00d767=6365\", \"48\u00d782=3936\"),\n    @(\"16\u00d725=400\", \"88\u00d757=5016\"),\n    @(\"43\u00d756=2408\", \"46\u00d714=644\"),\n    @(\"88\u00d731=2728\", \"62\u00d773=4526\"),\n    @(\"42\u00d797=4074\", \"92\u00d767=6164\"),\n    @(\"94\u00d749=4606\", \"19\u00d743=817\"),\n    @(\"38\u00d721=798\", \"90\u00d753=4770\"),\n    @(\"13\u00d780=1040\", \"64\u00d768=4352\"),\n    @(\"80\u00d751=4080\", \"29\u00d793=2697\"),\n    @(\"40\u00d716=640\", \"74\u00d745=3330\"),\n    @(\"45\u00d728=1260\", \"92\u00d745=4140\"),\n    @(\"58\u00d796=5568\", \"16\u00d712=192\"),\n    @(\"27\u00d744=1188\", \"56\u00d796=5376\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for '$oldText'\"\n    }\n}\n\n$d.Save()\n"}
